$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; existing B:L shift right to C:M.
# Capture the current column's width first so the newly inserted column
# keeps the same width as its neighbours.
$origWidth = $ws.Columns("B:B").ColumnWidth
$ws.Columns("B:B").Insert()
$ws.Columns("B:B").ColumnWidth = $origWidth

# New latest-price-check timestamp goes in the new column's header cell.
$ws.Range("B1").Value = "2025-12-21 14:19"

# The newly inserted column carries forward the previous "latest" price
# column's values (now shifted into column C), row by row.
$lastRow = 26
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Copy() | Out-Null
    $ws.Cells.Item($r, 2).PasteSpecial() | Out-Null
}

# Row 2 has no new sample this time around (matches the sheet's existing
# "missing data point" pattern), so its new column-B cell stays blank.
$ws.Range("B2").ClearContents()
